# Update scraped_at timestamps (column K) on the "snapshot" sheet
$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")

$updates = @(
    @{ Row = 2; Value = "2025-12-18T03:01:12.213741+00:00" },
    @{ Row = 3; Value = "2025-12-18T03:01:12.213775+00:00" },
    @{ Row = 4; Value = "2025-12-18T03:01:14.540583+00:00" },
    @{ Row = 5; Value = "2025-12-18T03:01:14.540614+00:00" },
    @{ Row = 6; Value = "2025-12-18T03:01:14.540631+00:00" },
    @{ Row = 7; Value = "2025-12-18T03:01:16.464099+00:00" },
    @{ Row = 8; Value = "2025-12-18T03:01:18.429264+00:00" },
    @{ Row = 9; Value = "2025-12-18T03:01:20.356173+00:00" },
    @{ Row = 10; Value = "2025-12-18T03:01:20.356202+00:00" },
    @{ Row = 11; Value = "2025-12-18T03:01:22.653884+00:00" },
    @{ Row = 12; Value = "2025-12-18T03:01:26.876766+00:00" },
    @{ Row = 13; Value = "2025-12-18T03:01:26.876796+00:00" },
    @{ Row = 14; Value = "2025-12-18T03:01:29.229775+00:00" },
    @{ Row = 15; Value = "2025-12-18T03:01:31.663909+00:00" },
    @{ Row = 16; Value = "2025-12-18T03:01:34.106076+00:00" },
    @{ Row = 17; Value = "2025-12-18T03:01:36.526272+00:00" },
    @{ Row = 18; Value = "2025-12-18T03:01:36.526300+00:00" },
    @{ Row = 19; Value = "2025-12-18T03:01:36.526318+00:00" },
    @{ Row = 20; Value = "2025-12-18T03:01:36.526333+00:00" },
    @{ Row = 21; Value = "2025-12-18T03:01:38.467760+00:00" },
    @{ Row = 22; Value = "2025-12-18T03:01:38.467790+00:00" },
    @{ Row = 23; Value = "2025-12-18T03:01:40.364269+00:00" },
    @{ Row = 24; Value = "2025-12-18T03:01:40.364297+00:00" },
    @{ Row = 25; Value = "2025-12-18T03:01:40.364314+00:00" },
    @{ Row = 26; Value = "2025-12-18T03:01:42.319743+00:00" },
    @{ Row = 27; Value = "2025-12-18T03:01:42.319773+00:00" },
    @{ Row = 28; Value = "2025-12-18T03:01:44.622837+00:00" },
    @{ Row = 29; Value = "2025-12-18T03:01:44.622870+00:00" },
    @{ Row = 30; Value = "2025-12-18T03:01:44.622888+00:00" },
    @{ Row = 31; Value = "2025-12-18T03:01:46.984438+00:00" },
    @{ Row = 32; Value = "2025-12-18T03:01:49.398458+00:00" },
    @{ Row = 33; Value = "2025-12-18T03:01:49.398484+00:00" },
    @{ Row = 34; Value = "2025-12-18T03:01:54.305801+00:00" },
    @{ Row = 35; Value = "2025-12-18T03:01:54.305843+00:00" },
    @{ Row = 36; Value = "2025-12-18T03:01:56.612556+00:00" },
    @{ Row = 37; Value = "2025-12-18T03:01:56.612583+00:00" }
)

foreach ($u in $updates) {
    $snapshot.Cells.Item($u.Row, 11).Value = $u.Value
}

# Remove all data rows (keep header) from the "returned" sheet
$returned = $wb.Worksheets.Item("returned")
$returned.Rows("2:9").Delete()

# Remove all data rows (keep header) from the "new_injured" sheet
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows("2:2").Delete()
